$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Positive Testdata" ---
$ws1 = $wb.Worksheets.Item("Positive Testdata")

$ws1.Range("A2").Value = "sk@g.com"
$ws1.Range("C2").Value = "Test@123"
$ws1.Range("D2").Value = "Test@123"
$ws1.Range("A3").Value = "anuj.lpu1@gmail.com"
$ws1.Range("B3").Value = 9009857868
$ws1.Range("C3").Value = "Test@123"
$ws1.Range("D3").Value = "Test@123"

# New hyperlink on A2 (preserve the existing non-hyperlink cell style)
$preservedStyle = $ws1.Range("A2").Style
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:sk@g.com")
$ws1.Range("A2").Style = $preservedStyle

$ws1.Columns.Item(2).ColumnWidth = 10

$ws1.Range("D3").Select()

# --- Sheet 2: "Negative Testdata" ---
$ws2 = $wb.Worksheets.Item("Negative Testdata")

$ws2.Range("C2").Value = "wk@1234"

$ws2.Columns.Item(4).ColumnWidth = 10.3

$ws2.Range("F6").Select()
